# Clean up non-public new spells
# - Adds 9 new "New" spells (Development Status: Playtest Ready / Release
#   Status: Not Released) to the bottom of the Spells sheet.
# - Leaves the active tab on "Subclasses" (it was left there after the
#   author finished editing "Spells").

$wb = $excel.ActiveWorkbook

$spells = $wb.Worksheets.Item("Spells")

# Columns: A=Name, B=New/Revised, C=Spell Level, D=School,
#          E=Artificer, F=Bard, G=Cleric, H=Druid, I=Paladin, J=Ranger,
#          K=Sorcerer, L=Warlock, M=Wizard,
#          N=Development Status, O=Release Status
$newSpellRows = @(
    @("Binding Chain",    "New", 1, "Conjuration",  "Yes", "No",  "No", "No",  "Yes", "Yes", "No",  "Yes", "No",  "Playtest Ready", "Not Released"),
    @("Molten Sphere",    "New", 6, "Conjuration",  "No",  "No",  "No", "Yes", "No",  "No",  "Yes", "No",  "No",  "Playtest Ready", "Not Released"),
    @("Mud Ball",         "New", 2, "Conjuration",  "No",  "No",  "No", "Yes", "No",  "Yes", "No",  "No",  "No",  "Playtest Ready", "Not Released"),
    @("Tranquility",      "New", 5, "Conjuration",  "No",  "No",  "No", "Yes", "No",  "No",  "No",  "No",  "No",  "Playtest Ready", "Not Released"),
    @("Frozen Tomb",      "New", 5, "Evocation",    "No",  "No",  "No", "Yes", "No",  "No",  "No",  "Yes", "No",  "Playtest Ready", "Not Released"),
    @("Water Whip",       "New", 1, "Conjuration",  "No",  "No",  "No", "Yes", "No",  "No",  "No",  "No",  "No",  "Playtest Ready", "Not Released"),
    @("Grasping Tide",    "New", 2, "Conjuration",  "No",  "No",  "No", "Yes", "No",  "No",  "No",  "No",  "No",  "Playtest Ready", "Not Released"),
    @("Stream of Flames", "New", 6, "Evocation",    "No",  "No",  "No", "No",  "No",  "No",  "Yes", "No",  "Yes", "Playtest Ready", "Not Released"),
    @("Leap Slam",        "New", 3, "Trasmutation", "No",  "No",  "No", "No",  "Yes", "No",  "No",  "No",  "No",  "Playtest Ready", "Not Released")
)

$columns = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O")

$startRow = 39
for ($i = 0; $i -lt $newSpellRows.Count; $i++) {
    $row = $startRow + $i
    $values = $newSpellRows[$i]
    for ($c = 0; $c -lt $columns.Count; $c++) {
        $spells.Range("$($columns[$c])$row").Value = $values[$c]
    }
}

# Leave "Subclasses" as the active/selected sheet (matches the saved
# workbook state after this edit).
$subclasses = $wb.Worksheets.Item("Subclasses")
$subclasses.Activate()
